# metro_budget -> "worked" cleanup pass:
#  - drop the helper "Avg_diff" column (K) and its AVERAGE() formulas
#  - drop the scratch INDEX/MATCH lookup helper cells (B:G) in the small
#    Departments table at the bottom of the sheet, keeping the department
#    names in column A
#  - remove the now-unused defined names
#  - tweak the data_dictionary description font to Arial
#  - leave data_dictionary as the active/selected sheet

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("metro_budget")
$ws2 = $wb.Worksheets.Item("data_dictionary")

# --- metro_budget sheet -----------------------------------------------

# Remove the Avg_diff column (K) values/formulas entirely.
$ws1.Columns.Item(11).ClearContents()

# Remove the INDEX/MATCH/AVERAGE scratch formulas in the little lookup
# table (rows 56-61), keeping the department-name column (A).
$ws1.Range("B56:G61").ClearContents()

# The defined names that only existed to support the removed formulas.
$wb.Names.Item("AllDepartments").Delete()
$wb.Names.Item("dataset").Delete()
$wb.Names.Item("Departments").Delete()
$wb.Names.Item("Headers").Delete()

# Selection left on the (now empty) column K by whoever did the cleanup.
$ws1.Range("K1:K1048576").Select()

# --- data_dictionary sheet ----------------------------------------------

# Description column font switched from Calibri to Arial.
$ws2.Range("B1:B10").Font.Name = "Arial"

# data_dictionary ends up as the active/selected sheet.
$ws2.Activate()
$ws2.Range("B14").Select()
